$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 29

$ws.Range("F6").Value = 28
$ws.Range("H6").Value = 35

$ws.Range("E10").Value = 35

$ws.Range("E11").Value = 22
$ws.Range("F11").Value = 17
$ws.Range("H11").Value = 18

$ws.Range("E12").Value = 32

$ws.Range("E14").Value = 39

$ws.Range("E15").Value = 105
$ws.Range("F15").Value = 52
$ws.Range("H15").Value = 63

$ws.Range("E16").Value = 324
$ws.Range("G16").Value = 88
$ws.Range("H16").Value = 182
